$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.157.58'
$ws.Range("E2").Value = '  -1.31%  '
$ws.Range("D3").Value = '1.802.44'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '223.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.553'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.44%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.37'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.287'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0719'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0927'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.82%  '
$ws.Range("D12").Value = '2.062.05'
$ws.Range("E12").Value = '  +0.57%  '
$ws.Range("D13").Value = '1.801.66'
$ws.Range("E13").Value = '  -0.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.86'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.10%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.632'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.44%  '
$ws.Range("D16").Value = '34.186.21'
$ws.Range("E16").Value = '  -1.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.23'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '68.54'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '248.16'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.15%  '
$ws.Range("D20").Value = '0.0₃0789'
$ws.Range("E20").Value = '  -0.25%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.00'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.94'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.95%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.12'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.12'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '159.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.60'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.113'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.10%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0526'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.95%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.21'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.70%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.52'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.11%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.53%  '
$ws.Range("D35").Value = '1.415.66'
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.651'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.45%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.06'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0187'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.56%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.948'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '80.47'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.37'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.74'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.96'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '107.68'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.17%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0494'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.15%  '
$ws.Range("D47").Value = '1.960.82'
$ws.Range("E47").Value = '  +0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '12.03'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.999'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("D51").Value = '0.0₆0125'
$ws.Range("E51").Value = '  +3.89%  '
